$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A33").Value = "Wat zijn jullie openingstijden?"
$ws.Range("B33").Value = "mailmind.test@zohomail.eu"
$ws.Range("C33").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$ws.Range("D33").Value = "Informatieaanvraag"
$ws.Range("E33").Value = "Beste [Naam],`nDank voor je interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 09:00 tot 18:00 uur. Op zaterdag zijn we geopend van 10:00 tot 16:00 uur. Op zondag zijn we gesloten. Mocht je verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$ws.Range("F33").Value = "2025-06-17 21:50:38"
$ws.Range("G33").Value = "Ja"

$ws.Rows.Item(33).AutoFit()

$dCfs = $ws.Range("D2:D32").FormatConditions
for ($i = 1; $i -le $dCfs.Count; $i++) {
    $dCfs.Item($i).ModifyAppliesToRange($ws.Range("D2:D33"))
}

$gCfs = $ws.Range("G2:G32").FormatConditions
for ($i = 1; $i -le $gCfs.Count; $i++) {
    $gCfs.Item($i).ModifyAppliesToRange($ws.Range("G2:G33"))
}

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 16
